$d = $word.ActiveDocument

# --- 1. Bold the first paragraph: "Vous devrez présenter une planification..." ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.Bold = 1
$r1.BoldBi = 1

# --- 2. Fix the last paragraph ("...doit être au pire le 31 mai)") so the
#         run split caused by the old _GoBack bookmark is merged back into a
#         single run containing the full, uninterrupted text. ---
$p12 = $d.Paragraphs.Item(12)
$r12 = $p12.Range
$oldText = "• Le projet doit être échelonné pour un rendu prévu le 31 mai (à vous de définir des dates de début, de fin et de jalons, mais la date de fin doit être au pire le 31 mai)"
$find12 = $r12.Find
$find12.ClearFormatting()
$find12.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2)

# --- 3. Bold the "Vous pouvez utiliser..." paragraph and split it where the
#         _GoBack bookmark used to sit, re-creating the bookmark there. ---
$p10 = $d.Paragraphs.Item(10)
$r10 = $p10.Range
$r10.Bold = 1
$r10.BoldBi = 1

# Locate the end of the word "planification" inside the paragraph - that is
# exactly where the bookmark sits, right before " : ".
$locate = $p10.Range.Duplicate
$findWord = $locate.Find
$findWord.ClearFormatting()
$findWord.Text = "planification"
$findWord.Forward = $true
$findWord.Wrap = 0
$findWord.Execute() | Out-Null
$splitOffset = $locate.End

$bmRange = $d.Range($splitOffset, $splitOffset)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Re-stamp the first half of the split so the serializer doesn't keep a
# stray (unneeded) xml:space="preserve" on a run with no leading/trailing
# whitespace.
$firstHalf = $d.Range($p10.Range.Start, $splitOffset)
$findFirstHalf = $firstHalf.Find
$findFirstHalf.ClearFormatting()
$firstHalfText = "Vous pouvez utiliser les informations suivantes dans votre planification"
$findFirstHalf.Execute($firstHalfText, $true, $false, $false, $false, $false, $true, 1, $false, $firstHalfText, 2)
